# Update cryptos list values (price/volume) per upstream data refresh.
# Also corrects the BabyDogeCoin/ApeXProtocol row order (rows 49-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.500.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.932.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000355"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.568.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.935.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.546.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  +19.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "722.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0913"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "61.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.02%  "
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +18.39%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.73%  "
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.23%  "
$ws.Range("E44").Value = "  +6.08%  "
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0349"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +40.77%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "146.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
